# Update Name of Algo
# Applies updated numeric results for the RandomForest imputation output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = -10.9086
$ws.Range("A3").Value  = -21.46100000000003
$ws.Range("C5").Value  = -14.49230000000001
$ws.Range("D5").Value  = -8.720499999999992
$ws.Range("E7").Value  = 12.1538
$ws.Range("D9").Value  = -8.687400000000002
$ws.Range("D11").Value = -8.182600000000001
$ws.Range("E11").Value = 13.4199
$ws.Range("A14").Value = -20.61979999999999
$ws.Range("A16").Value = -20.43189999999999
$ws.Range("C16").Value = -11.49889999999999
$ws.Range("D17").Value = -8.582100000000002
$ws.Range("E19").Value = 13.26859999999999
$ws.Range("A21").Value = -21.48160000000001
$ws.Range("D21").Value = -7.493100000000007
$ws.Range("E21").Value = 13.30600000000001
$ws.Range("A23").Value = -21.37300000000003
$ws.Range("A25").Value = -22.49320000000004
